# Update the "2. Data reporter" block (Organization / Contact person / email /
# phone / website) with the new National Statistical Committee contact
# information, then leave the active selection on the Organization cell (B6)
# to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value  = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value  = "Kalymbetova Yryskan"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

$ws.Range("B6").Select()
